{"js": "// Add two new rows to the end of the (only) table in the document body,\n// mirroring the \"datetimeField\" and \"objectLevel1\" schema rows added by\n// the diff.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\ntable.addRows(Word.InsertLocation.end, 2, [\n  [\"datetimeField\", \"Datetime\", \"datetime\", \"0..1\", \"Datetime\", \"\"],\n  [\"objectLevel1\", \"Level 1 Object\", \"object\", \"0..1\", \"Object at data level 1\", \"\"],\n]);\n\nawait context.sync();\n", "ps1": "# Add two new rows to the end of the (only) table in the document,\n# mirroring the \"datetimeField\" and \"objectLevel1\" schema rows added by\n# the diff.\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$row1 = $t.Rows.Add()\n$row1.Cells(1).Range.Text = \"datetimeField\"\n$row1.Cells(2).Range.Text = \"Datetime\"\n$row1.Cells(3).Range.Text = \"datetime\"\n$row1.Cells(4).Range.Text = \"0..1\"\n$row1.Cells(5).Range.Text = \"Datetime\"\n\n$row2 = $t.Rows.Add()\n$row2.Cells(1).Range.Text = \"objectLevel1\"\n$row2.Cells(2).Range.Text = \"Level 1 Object\"\n$row2.Cells(3).Range.Text = \"object\"\n$row2.Cells(4).Range.Text = \"0..1\"\n$row2.Cells(5).Range.Text = \"Object at data level 1\"\n"}
